$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B18").Value = "[Humberto-Desenho tecnico mecanico, Weslei-Desenho tecnico mecanico]"
$ws.Range("C18").Value = "[Suzanny-Metalografia, Mayra-Tec. Mat. Não Metal., Suzanny-Trat. Termicos, Victor-Ajustagem]"
$ws.Range("D18").Value = "[Aselmo-Manut. Mot. End., Weslei-Metrologia 1, Andre B.-Elet. Dig. Bas., Valmir-Caldeiraria]"
$ws.Range("E18").Value = "Gilberto-M.T.R"
$ws.Range("F18").Value = "[Sandro-Comandos Eletricos, Victor-Ajustagem, Gisele-E. D. N. D., Valmir-Caldeiraria]"

$ws.Range("B19").Value = "[Humberto-Desenho tecnico mecanico, Weslei-Desenho tecnico mecanico]"
$ws.Range("C19").Value = "[Suzanny-Metalografia, Mayra-Tec. Mat. Não Metal., Weslei-Metrologia 1, Anderson-Tornearia]"
$ws.Range("D19").Value = "[Rachel-T.M. Metalicos, Gisele-E. D. N. D., Andre B.-Elet. Dig. Bas., Weslei-Metrologia 1]"
$ws.Range("E19").Value = "Humberto-M.T"
$ws.Range("F19").Value = "[Sandro-Comandos Eletricos, Aselmo-Manut. Mot. End., Andre B.-Elet. Dig. Bas., Valmir-Caldeiraria]"

$ws.Range("B20").Value = "[Humberto-Desenho tecnico mecanico, Weslei-Desenho tecnico mecanico]"
$ws.Range("C20").Value = "[Suzanny-Trat. Termicos, Suzanny-Metalografia, Anderson-Tornearia, Victor-Ajustagem]"
$ws.Range("D20").Value = "[Rachel-T.M. Metalicos, Anderson-Tornearia, Mayra-Tec. Mat. Não Metal., Andre B.-Elet. Dig. Bas.]"
$ws.Range("E20").Value = "Humberto-M.T"
$ws.Range("F20").Value = "[Sandro-Comandos Eletricos, Aselmo-Manut. Mot. End., Gisele-E. D. N. D., Rachel-T.M. Metalicos]"

$ws.Range("B21").Value = "[Victor-Ajustagem, Mayra-Tec. Mat. Não Metal., Anderson-Tornearia, Suzanny-Trat. Termicos]"
$ws.Range("C21").Value = "Gilberto-M.T.R"
$ws.Range("D21").Value = "[Suzanny-Metalografia, Suzanny-Trat. Termicos, Sandro-Comandos Eletricos, Weslei-Metrologia 1]"
$ws.Range("E21").Value = "Gilberto-M.T.R"
$ws.Range("F21").Value = "[Rachel-T.M. Metalicos, Aselmo-Manut. Mot. End., Gisele-E. D. N. D., Valmir-Caldeiraria]"
